# ordenanzas/0099.docx - reformat headings (bold/underline/keepNext),
# split VISTO/CONSIDERANDO lead-ins into their own paragraphs, tighten
# the long run of spaces before "(barrio las orquidias)", add indentation
# to the "EL CONCEJO..." sanction line, and start page numbering at 65
# with a footer on the section.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 1: "Yerba Buena, 21 de Noviembre de 1984"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Format.KeepWithNext = $true
$p1.Format.SpaceAfter = 12

# ---------------------------------------------------------------------
# Paragraph 2: "ORDENANZA Nº 99"
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Format.KeepWithNext = $true
$p2.Format.SpaceBefore = 12
$p2.Format.SpaceAfter = 18
$p2.Range.Font.Bold = $true

# ---------------------------------------------------------------------
# Paragraph 3: "VISTO: Que en fecha ..." -> split into "VISTO: " heading
# paragraph and a body paragraph that starts with a single space.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$f3 = $r3.Find
$f3.Execute("VISTO: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Collapse(0)
$r3.InsertParagraphAfter()

$visto = $d.Paragraphs(3)
$visto.Format.KeepWithNext = $true
$visto.Format.SpaceBefore = 12
$visto.Format.SpaceAfter = 6
$visto.Range.Font.Bold = $true

$body3 = $d.Paragraphs(4)
$body3.Format.KeepWithNext = $true
$body3.Format.SpaceAfter = 6
$b3r = $body3.Range
$b3r.Collapse(1)
$b3r.InsertBefore(" ")
$b3space = $body3.Range.Characters(1)
$b3space.Font.Name = "Times New Roman"
$b3space.Font.NameBi = "Times New Roman"
$b3space.Font.Size = 12

# Tighten the run of 20 spaces before "(barrio las orquidias)" down to 1.
$find = $d.Content.Find
$find.Execute("                    (", $true, $false, $false, $false, $false, $true, 1, $false, " (", 2)

# ---------------------------------------------------------------------
# Paragraph 5: "CONSIDERANDO: Que, como repesentates ..." -> split the
# same way as VISTO.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range
$f5 = $r5.Find
$f5.Execute("CONSIDERANDO: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r5.Collapse(0)
$r5.InsertParagraphAfter()

$consid = $d.Paragraphs(5)
$consid.Format.KeepWithNext = $true
$consid.Format.SpaceBefore = 12
$consid.Format.SpaceAfter = 6
$consid.Range.Font.Bold = $true

$body5 = $d.Paragraphs(6)
$body5.Format.KeepWithNext = $true
$body5.Format.SpaceAfter = 6
$b5r = $body5.Range
$b5r.Collapse(1)
$b5r.InsertBefore(" ")
$b5space = $body5.Range.Characters(1)
$b5space.Font.Name = "Times New Roman"
$b5space.Font.NameBi = "Times New Roman"
$b5space.Font.Size = 12

# ---------------------------------------------------------------------
# Paragraph 7: "Por ello:"
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$p7.Format.KeepWithNext = $true
$p7.Format.SpaceAfter = 6

# ---------------------------------------------------------------------
# Paragraph 8: "EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA"
# ---------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$p8.Format.KeepWithNext = $true
$p8.Format.SpaceBefore = 18
$p8.Format.SpaceAfter = 18
$p8.Format.LeftIndent = 99.2
$p8.Format.RightIndent = 99.2
$p8.Range.Font.Bold = $true

# ---------------------------------------------------------------------
# Paragraphs 9-12: "ARTICULO PRIMERO/SEGUNDO/TERCERO/CUARTO: "
# Each gets keepNext + SpaceAfter=6 (no SpaceBefore), loses any
# justification, and its label becomes underlined with the trailing
# space split into its own (non-underlined) run.
# ---------------------------------------------------------------------
$labels = @("ARTICULO PRIMERO: ", "ARTICULO SEGUNDO: ", "ARTICULO TERCERO: ", "ARTICULO CUARTO: ")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $para = $d.Paragraphs(9 + $i)
    $para.Format.Alignment = 0
    $para.Format.KeepWithNext = $true
    $para.Format.SpaceAfter = 6

    $pr = $para.Range
    $pf = $pr.Find
    $pf.Execute($labels[$i], $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $label = $d.Range($pr.Start, $pr.End - 1)
    $label.Font.Underline = 1
}

# ---------------------------------------------------------------------
# Section: footer + page numbering starting at 65.
# ---------------------------------------------------------------------
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$ftr.PageNumbers.Add()
$ftr.PageNumbers.StartingNumber = 65

Write-Output "done"
